# Applies the "Made final edits and ran tests." commit:
#  - fills in the Win/Loss trial grid (H3:AC3, E4:AC6)
#  - adds AD/AE COUNTIF summary formulas for every row
#  - adds a new "Overall results" column (AF) with a merged header and a
#    merged, wrapped commentary cell next to the first trial row
#  - widens column AF and updates the sheet selection accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlLeft   = -4131

# ---------------------------------------------------------------------------
# 1. Win (W) / Loss (L) grid for the four trial rows (columns E..AC)
# ---------------------------------------------------------------------------
$row3 = @("L","L","L","L","L","L","L","L","L","L","L","W","L","W","L","L","L","L","L","L","L","L","L","L","L")
$row4 = @("L","L","W","W","W","W","W","W","W","W","W","W","L","W","W","W","W","W","W","W","W","L","W","W","W")
$row5 = @("L","W","W","W","L","W","L","W","W","L","L","L","L","W","L","L","L","L","L","L","L","L","L","L","L")
$row6 = @("L","L","L","L","L","L","L","L","L","L","L","L","L","L","L","W","L","L","L","L","L","L","W","L","L")

$rows = @{ 3 = $row3; 4 = $row4; 5 = $row5; 6 = $row6 }

foreach ($r in 3..6) {
    $data = $rows[$r]
    for ($i = 0; $i -lt $data.Length; $i++) {
        # column E is index 5
        $ws.Cells.Item($r, 5 + $i).Value = $data[$i]
    }

    # ------------------------------------------------------------------
    # 2. Win / Loss tallies
    # ------------------------------------------------------------------
    $ws.Cells.Item($r, 30).Formula = '=COUNTIF($E${0}:$AC${0}, "W")' -f $r
    $ws.Cells.Item($r, 31).Formula = '=COUNTIF($E${0}:$AC${0}, "L")' -f $r
}

# ---------------------------------------------------------------------------
# 3. New column AF: "Overall results" header + commentary
# ---------------------------------------------------------------------------
$ws.Range("AF1").Value = "Overall results"
$ws.Range("AF3").Value = "If the player selects the Truck or Light vehicle, it is possible to win more often than when all vehicles are controlled by the algorithm. However, this is not shown by the results when controlling the Heavy vehicle. If a future version is released, it may be useful to make the Heavy and Light vehicles easier to control."

# Header + body merges
$ws.Range("AF1:AF2").Merge()
$ws.Range("AF3:AF6").Merge()

# Style: centered horizontally, wrapped text (new cellXfs entry)
$ws.Range("AF1:AF6").HorizontalAlignment = $xlCenter
$ws.Range("AF1:AF6").WrapText = $true

# Column width to fit the commentary text
$ws.Range("AF1").ColumnWidth = 65.83

# ---------------------------------------------------------------------------
# 4. Update the sheet selection / view to focus on the new column
# ---------------------------------------------------------------------------
$ws.Range("D1").Select()
$ws.Range("AF3:AF6").Select()
